## new results in table1 and 2
# Fill in the "original_subproblem" sheet (table1) with the new run's
# results, and move the visible selection from "original_validCuts"
# (table2) over to "original_subproblem" which becomes the active tab.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("original_subproblem")
$ws2 = $wb.Worksheets.Item("original_validCuts")

# --- original_subproblem (table1): new run results -----------------------
$ws1.Range("A1").Value = "Budget :"
$ws1.Range("B1").Value = 30000000
$ws1.Range("C1").Value = "Tornado Length:"
$ws1.Range("D1").Value = 5

$ws1.Range("A2").Value = "Best Bound:"
$ws1.Range("B2").Value = 0

$ws1.Range("A3").Value = "Best Objective:"
$ws1.Range("B3").Value = 103

$ws1.Range("A4").Value = "Gap:"
$ws1.Range("B4").Value = "-"

$ws1.Range("A5").Value = "CCG Run time:"
$ws1.Range("B5").Value = 3600.26

$ws1.Range("A6").Value = "CCG Iteration:"
$ws1.Range("B6").Value = 1

$ws1.Range("A7").Value = "Subproblem Run time:"
$ws1.Range("B7").Value = 3600.22

$ws1.Range("A8").Value = "Number of Subproblem Callbacks:"
$ws1.Range("B8").Value = 5

$ws1.Range("A9").Value = "Subproblem Callbacks Run Time:"
$ws1.Range("B9").Value = 0.21

$ws1.Range("A10").Value = "Number of Uncertainty Set Check Call:"
$ws1.Range("B10").Value = 0

$ws1.Range("A11").Value = "Uncertainty Set Check Run Time:"
$ws1.Range("B11").Value = 0

# --- view / selection state ----------------------------------------------
# table2 (original_validCuts) keeps its whole new-results block selected ...
[void]$ws2.Activate()
[void]$ws2.Range("A1:D11").Select()

# ... but table1 (original_subproblem) ends up the active tab with B8 selected
[void]$ws1.Activate()
[void]$ws1.Range("B8").Select()
